# Mise à jour de l'application
# Adds a new attendance column (AH) for the session dated 2025-08-26
# (serial 45895), mirroring the formatting of the existing AG column
# and filling in each player's attendance status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: date for the new session (values are written before the
# formatting pass below so the dependent COUNTA/COUNTIF formulas pick
# up the new cells and recalculate correctly).
$ws.Cells.Item(1, 34).Value = 45895

# Attendance values for each player row (2-26) and the last row (27)
$presentRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,18,19,20,21,22,23,24,25,26)
foreach ($r in $presentRows) {
    $ws.Cells.Item($r, 34).Value = "P"
}

# Row 17 (Maé Clavel) was injured ("B" = Blessure) that day
$ws.Cells.Item(17, 34).Value = "B"

# Row 27 (Sofiane Belle) was off that day ("RH")
$ws.Cells.Item(27, 34).Value = "RH"

# Copy the formatting of the last existing day column (AG) onto the new
# column (AH) for every row so styles (date header style / data style)
# match exactly.
$ws.Range("AG1:AG27").Copy()
$ws.Range("AH1:AH27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to match the edited workbook's last cursor spot
$ws.Range("AJ25").Select()
